# Apply the beta/xi -> beta_A / beta_P notation change in the heroin
# model schematic diagram.
#
#   beta(1-xi)SA  ->  beta_A SA    (TextBox 57 / shape "TextBox 57")
#   beta*xi*SP    ->  beta_P SP    (TextBox 38 / shape "TextBox 38")
#
# In both cases the letter that becomes a subscript ("A" or "P") is
# rendered with a lowered baseline (baseline = -25%), and the textbox is
# nudged slightly to the right to keep it visually centered over the
# arrow it annotates.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# Shape 10 "TextBox 57": "beta(1-xi)SA"  ->  "beta" + "A " (sub) + "SA"
# ---------------------------------------------------------------------
$shape1 = $s.Shapes.Item(10)

# Re-center the textbox horizontally (position changed in the diff).
$shape1.Left = 2875824 / 12700

$tr1 = $shape1.TextFrame.TextRange
# Go through a throwaway value first so the final assignment doesn't
# keep stray formatting (e.g. the spell-check "err" flag) bleeding in
# from characters that used to occupy those positions.
$tr1.Text = "X"
$tr1.Text = [char]0x03B2 + "A SA"

$sub1 = $tr1.Characters(2, 2)
$sub1.Font.Subscript = $true
$sub1.Font.Name = "Times New Roman"

# ---------------------------------------------------------------------
# Shape 20 "TextBox 38": "beta*xi*SP"  ->  "beta" + "P" (sub) + "SP"
# ---------------------------------------------------------------------
$shape2 = $s.Shapes.Item(20)

# Re-center the textbox horizontally (position changed in the diff).
$shape2.Left = 3033659 / 12700

$tr2 = $shape2.TextFrame.TextRange
$tr2.Text = "X"
$tr2.Text = [char]0x03B2 + "PSP"

$sub2 = $tr2.Characters(2, 1)
$sub2.Font.Subscript = $true
$sub2.Font.Name = "Times New Roman"
